$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (George Thompson / Experienced referral) ---
$ws.Range("G2").Value = "George.Thompson1a973@gmail.com"
$ws.Range("H2").Value = 3333333333

# --- Add a new row 3 (second "refer a friend" submission - a Fresher) ---
# Start from a copy of row 2 so formatting/styles (borders, hyperlink style, etc.) match.
$ws.Range("A2:P2").Copy($ws.Range("A3:P3"))

# Fields that stay the same as row 2
# (A3 UserName, B3 Password, C3 FirstName, E3 LastName, F3 Gender, M3 Resume,
#  O3 Country, P3 SearchKeyword all copied already from row 2)

# Fields that differ for this candidate
$ws.Range("G3").Value = "George.Thompson2@gmail.com8888888888"
$ws.Range("H3").Value = 8788888888
$ws.Range("I3").Value = "Fresher"
$ws.Range("N3").Value = "Available for immediate joining"

# Fields that are blank for this (fresher) candidate
$ws.Range("D3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()

# Hyperlink the new email address cell, same as row 2's G column
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:George.Thompson2@gmail.com8888888888") | Out-Null

# Re-apply row2's hyperlink cell style to G3 (Hyperlinks.Add resets it to a plain
# hyperlink style without the row's border) so G3 matches G2's look.
$ws.Range("G2").Copy($ws.Range("G3"))
$ws.Range("G3").Value = "George.Thompson2@gmail.com8888888888"

# --- Selection moves to G2 ---
$ws.Range("G2").Select() | Out-Null
